# Update the workbook to reflect data through 2022-11-20
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet / tab label
$ws.Name = "Through 2022-11-20"

# Update header label in I1 (shared string "2022 (through 11-14)" -> "2022 (through 11-20)")
$ws.Range("I1").Value = "2022 (through 11-20)"

# Update the monthly figures for the current year column (I)
$ws.Range("I10").Value = 143    # September
$ws.Range("I12").Value = 76     # November
$ws.Range("I14").Value = 1473   # Total
